# Natmi following Dr Hou advice
# Update LR-pair results table (Hspg2 -> Col13a1) with the recomputed
# natmi statistics (now using 3 replicate ligand/receptor expressing-cell
# counts and weighted cluster combinations, expanding the table from
# 3 data rows to 6 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Hspg2"
$ws.Cells.Item(2, 3).Value = "Col13a1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 133.9646796666667
$ws.Cells.Item(2, 8).Value = 401.894039
$ws.Cells.Item(2, 9).Value = 0.2795129415517746
$ws.Cells.Item(2, 10).Value = 0.2795129415517745
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.3326186666666667
$ws.Cells.Item(2, 14).Value = 0.9978560000000001
$ws.Cells.Item(2, 15).Value = 0.3213922220876632
$ws.Cells.Item(2, 16).Value = 0.3213922220876632
$ws.Cells.Item(2, 17).Value = 44.55915313115378
$ws.Cells.Item(2, 18).Value = 401.0323781803841
$ws.Cells.Item(2, 19).Value = 0.08983328538758398
$ws.Cells.Item(2, 20).Value = 0.08983328538758395
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Hspg2"
$ws.Cells.Item(3, 3).Value = "Col13a1"
$ws.Cells.Item(3, 4).Value = "sCs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 133.9646796666667
$ws.Cells.Item(3, 8).Value = 401.894039
$ws.Cells.Item(3, 9).Value = 0.2795129415517746
$ws.Cells.Item(3, 10).Value = 0.2795129415517745
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.702312
$ws.Cells.Item(3, 14).Value = 2.106936
$ws.Cells.Item(3, 15).Value = 0.6786077779123368
$ws.Cells.Item(3, 16).Value = 0.6786077779123368
$ws.Cells.Item(3, 17).Value = 94.08500210605602
$ws.Cells.Item(3, 18).Value = 846.7650189545041
$ws.Cells.Item(3, 19).Value = 0.1896796561641907
$ws.Cells.Item(3, 20).Value = 0.1896796561641906
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Hspg2"
$ws.Cells.Item(4, 3).Value = "Col13a1"
$ws.Cells.Item(4, 4).Value = "ECs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 276.4348856666666
$ws.Cells.Item(4, 8).Value = 829.3046569999999
$ws.Cells.Item(4, 9).Value = 0.5767723868147629
$ws.Cells.Item(4, 10).Value = 0.5767723868147629
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.3326186666666667
$ws.Cells.Item(4, 14).Value = 0.9978560000000001
$ws.Cells.Item(4, 15).Value = 0.3213922220876632
$ws.Cells.Item(4, 16).Value = 0.3213922220876632
$ws.Cells.Item(4, 17).Value = 91.9474030905991
$ws.Cells.Item(4, 18).Value = 827.526627815392
$ws.Cells.Item(4, 19).Value = 0.1853701590372019
$ws.Cells.Item(4, 20).Value = 0.1853701590372019
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Hspg2"
$ws.Cells.Item(5, 3).Value = "Col13a1"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 276.4348856666666
$ws.Cells.Item(5, 8).Value = 829.3046569999999
$ws.Cells.Item(5, 9).Value = 0.5767723868147629
$ws.Cells.Item(5, 10).Value = 0.5767723868147629
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.702312
$ws.Cells.Item(5, 14).Value = 2.106936
$ws.Cells.Item(5, 15).Value = 0.6786077779123368
$ws.Cells.Item(5, 16).Value = 0.6786077779123368
$ws.Cells.Item(5, 17).Value = 194.143537422328
$ws.Cells.Item(5, 18).Value = 1747.291836800952
$ws.Cells.Item(5, 19).Value = 0.3914022277775611
$ws.Cells.Item(5, 20).Value = 0.3914022277775611
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Hspg2"
$ws.Cells.Item(6, 3).Value = "Col13a1"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 68.87942233333334
$ws.Cells.Item(6, 8).Value = 206.638267
$ws.Cells.Item(6, 9).Value = 0.1437146716334625
$ws.Cells.Item(6, 10).Value = 0.1437146716334625
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.3326186666666667
$ws.Cells.Item(6, 14).Value = 0.9978560000000001
$ws.Cells.Item(6, 15).Value = 0.3213922220876632
$ws.Cells.Item(6, 16).Value = 0.3213922220876632
$ws.Cells.Item(6, 17).Value = 22.91058161728356
$ws.Cells.Item(6, 18).Value = 206.195234555552
$ws.Cells.Item(6, 19).Value = 0.04618877766287739
$ws.Cells.Item(6, 20).Value = 0.04618877766287738
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Hspg2"
$ws.Cells.Item(7, 3).Value = "Col13a1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 68.87942233333334
$ws.Cells.Item(7, 8).Value = 206.638267
$ws.Cells.Item(7, 9).Value = 0.1437146716334625
$ws.Cells.Item(7, 10).Value = 0.1437146716334625
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.702312
$ws.Cells.Item(7, 14).Value = 2.106936
$ws.Cells.Item(7, 15).Value = 0.6786077779123368
$ws.Cells.Item(7, 16).Value = 0.6786077779123368
$ws.Cells.Item(7, 17).Value = 48.374844857768
$ws.Cells.Item(7, 18).Value = 435.3736037199121
$ws.Cells.Item(7, 19).Value = 0.09752589397058517
$ws.Cells.Item(7, 20).Value = 0.09752589397058516
